$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.67358826315247
$ws.Range("D2").Value = 8.289104047263637
$ws.Range("E2").Value = 14.29759936453496
$ws.Range("F2").Value = 40.38405106556283
$ws.Range("G2").Value = 3.701267103280734
$ws.Range("J2").Value = 11.01246239226328
$ws.Range("K2").Value = 18.97621323937035
$ws.Range("L2").Value = 9.726146552788437
$ws.Range("N2").Value = 18.59642389487522
$ws.Range("O2").Value = 30.97852337568994
$ws.Range("C3").Value = 13.63073567210559
$ws.Range("D3").Value = 8.263786884984988
$ws.Range("E3").Value = 14.29744455838229
$ws.Range("F3").Value = 40.45339490845652
$ws.Range("G3").Value = 3.703764472672328
$ws.Range("J3").Value = 11.03460911040082
$ws.Range("K3").Value = 18.62835538780363
$ws.Range("L3").Value = 9.739788378274355
$ws.Range("N3").Value = 18.64288455990328
$ws.Range("O3").Value = 31.06384755377812
$ws.Range("C4").Value = 13.6071832990912
$ws.Range("D4").Value = 8.249241341407334
$ws.Range("E4").Value = 14.29952732620293
$ws.Range("F4").Value = 40.50552379822315
$ws.Range("G4").Value = 3.705379283747629
$ws.Range("J4").Value = 11.04947621960687
$ws.Range("K4").Value = 18.41404032335869
$ws.Range("L4").Value = 9.748995781158044
$ws.Range("N4").Value = 18.67332379251888
$ws.Range("O4").Value = 31.12311698055277
$ws.Range("C5").Value = 13.59828632267558
$ws.Range("D5").Value = 8.243569005079978
$ws.Range("E5").Value = 14.30092493078742
$ws.Range("F5").Value = 40.5291628740659
$ws.Range("G5").Value = 3.706057872465073
$ws.Range("J5").Value = 11.0558540658845
$ws.Range("K5").Value = 18.32663210825692
$ws.Range("L5").Value = 9.752957326457066
$ws.Range("N5").Value = 18.6862097565437
$ws.Range("O5").Value = 31.14899517283063
$ws.Range("C6").Value = 13.59685150236846
$ws.Range("D6").Value = 8.242642631951226
$ws.Range("E6").Value = 14.30119018388234
$ws.Range("F6").Value = 40.53323266501042
$ws.Range("G6").Value = 3.706171794155473
$ws.Range("J6").Value = 11.05693240061442
$ws.Range("K6").Value = 18.31211688838039
$ws.Range("L6").Value = 9.753627800926964
$ws.Range("N6").Value = 18.68837858289011
$ws.Range("O6").Value = 31.15339629839092
$ws.Range("C7").Value = 13.60706046516416
$ws.Range("D7").Value = 8.249163804560594
$ws.Range("E7").Value = 14.29954395098973
$ws.Range("F7").Value = 40.5058329093181
$ws.Range("G7").Value = 3.705388352177326
$ws.Range("J7").Value = 11.04956093993566
$ws.Range("K7").Value = 18.41286165438265
$ws.Range("L7").Value = 9.749048359311011
$ws.Range("N7").Value = 18.6734956254477
$ws.Range("O7").Value = 31.12345900307572
$ws.Range("C8").Value = 13.65824439424625
$ws.Range("D8").Value = 8.280169632018088
$ws.Range("E8").Value = 14.29709477215297
$ws.Range("F8").Value = 40.40597521175071
$ws.Range("G8").Value = 3.702111338916336
$ws.Range("J8").Value = 11.01983539197326
$ws.Range("K8").Value = 18.85648855027079
$ws.Range("L8").Value = 9.730677980836118
$ws.Range("N8").Value = 18.61204718840844
$ws.Range("O8").Value = 31.00651248956358
$ws.Range("C9").Value = 13.78017158702268
$ws.Range("D9").Value = 8.348725397278269
$ws.Range("E9").Value = 14.30951159652516
$ws.Range("F9").Value = 40.28616645186906
$ws.Range("G9").Value = 3.696328055412074
$ws.Range("J9").Value = 10.97159909399153
$ws.Range("K9").Value = 19.71594459172837
$ws.Range("L9").Value = 9.701230751768215
$ws.Range("N9").Value = 18.50667998016993
$ws.Range("O9").Value = 30.83196895189689
$ws.Range("C10").Value = 13.88239416112919
$ws.Range("D10").Value = 8.40357336993042
$ws.Range("E10").Value = 14.32904139564084
$ws.Range("F10").Value = 40.24474690853583
$ws.Range("G10").Value = 3.692466726659435
$ws.Range("J10").Value = 10.9422723073727
$ws.Range("K10").Value = 20.3349151358831
$ws.Range("L10").Value = 9.683580210497825
$ws.Range("N10").Value = 18.43843880759392
$ws.Range("O10").Value = 30.73739041282923
$ws.Range("C11").Value = 13.9315228647973
$ws.Range("D11").Value = 8.429441506685755
$ws.Range("E11").Value = 14.34016247333653
$ws.Range("F11").Value = 40.23606364546801
$ws.Range("G11").Value = 3.690793368035108
$ws.Range("J11").Value = 10.93025421092351
$ws.Range("K11").Value = 20.61258793983411
$ws.Range("L11").Value = 9.676410270997026
$ws.Range("N11").Value = 18.40937484329884
$ws.Range("O11").Value = 30.70172267864775
$ws.Range("C12").Value = 13.95049325150053
$ws.Range("D12").Value = 8.439363962761004
$ws.Range("E12").Value = 14.34469302957848
$ws.Range("F12").Value = 40.23423804157484
$ws.Range("G12").Value = 3.690171602203706
$ws.Range("J12").Value = 10.92589317127829
$ws.Range("K12").Value = 20.71708050054399
$ws.Range("L12").Value = 9.673818329390745
$ws.Range("N12").Value = 18.39865289741028
$ws.Range("O12").Value = 30.6892776800666
$ws.Range("C13").Value = 13.94639153011275
$ws.Range("D13").Value = 8.437221431116818
$ws.Range("E13").Value = 14.34370313539469
$ws.Range("F13").Value = 40.2345661522551
$ws.Range("G13").Value = 3.690304982387393
$ws.Range("J13").Value = 10.92682395531111
$ws.Range("K13").Value = 20.69460679982938
$ws.Range("L13").Value = 9.67437107959101
$ws.Range("N13").Value = 18.40094944215613
$ws.Range("O13").Value = 30.69191066225563
$ws.Range("C14").Value = 13.93307629362988
$ws.Range("D14").Value = 8.430255323035286
$ws.Range("E14").Value = 14.34052882135024
$ws.Range("F14").Value = 40.23588413332364
$ws.Range("G14").Value = 3.690741976872415
$ws.Range("J14").Value = 10.92989162018156
$ws.Range("K14").Value = 20.62119828217447
$ws.Range("L14").Value = 9.676194564530977
$ws.Range("N14").Value = 18.4084870554154
$ws.Range("O14").Value = 30.70067752486339
$ws.Range("C15").Value = 13.92496769733719
$ws.Range("D15").Value = 8.426004731244481
$ws.Range("E15").Value = 14.33862596226517
$ws.Range("F15").Value = 40.23688193940313
$ws.Range("G15").Value = 3.691011196341382
$ws.Range("J15").Value = 10.93179538292256
$ws.Range("K15").Value = 20.57614526260617
$ws.Range("L15").Value = 9.677327528112974
$ws.Range("N15").Value = 18.41314101984465
$ws.Range("O15").Value = 30.70618583816649
$ws.Range("C16").Value = 13.87923539133472
$ws.Range("D16").Value = 8.401900849731755
$ws.Range("E16").Value = 14.32835941714247
$ws.Range("F16").Value = 40.24551896985096
$ws.Range("G16").Value = 3.692577753022708
$ws.Range("J16").Value = 10.94308431399079
$ws.Range("K16").Value = 20.31668182295452
$ws.Range("L16").Value = 9.684066040130334
$ws.Range("N16").Value = 18.44037797632142
$ws.Range("O16").Value = 30.73986974223978
$ws.Range("C17").Value = 13.8518450773513
$ws.Range("D17").Value = 8.387345261186759
$ws.Range("E17").Value = 14.32263249091844
$ws.Range("F17").Value = 40.25342093161373
$ws.Range("G17").Value = 3.693560044518064
$ws.Range("J17").Value = 10.95034832375105
$ws.Range("K17").Value = 20.15644108434562
$ws.Range("L17").Value = 9.688419710320643
$ws.Range("N17").Value = 18.4575934241142
$ws.Range("O17").Value = 30.762420671254
$ws.Range("C18").Value = 13.83633889124793
$ws.Range("D18").Value = 8.379059953774458
$ws.Range("E18").Value = 14.31954921755663
$ws.Range("F18").Value = 40.2589220683708
$ws.Range("G18").Value = 3.694132865519913
$ws.Range("J18").Value = 10.95465091514766
$ws.Range("K18").Value = 20.06391314064859
$ws.Range("L18").Value = 9.69100474055338
$ws.Range("N18").Value = 18.46768163294011
$ws.Range("O18").Value = 30.77608359585642
$ws.Range("C19").Value = 13.83113169365063
$ws.Range("D19").Value = 8.376269734069965
$ws.Range("E19").Value = 14.31854153272776
$ws.Range("F19").Value = 40.26094880790968
$ws.Range("G19").Value = 3.694328160085543
$ws.Range("J19").Value = 10.956129093925
$ws.Range("K19").Value = 20.03252563325582
$ws.Range("L19").Value = 9.691893897168056
$ws.Range("N19").Value = 18.47112935182809
$ws.Range("O19").Value = 30.78082839873802
$ws.Range("C20").Value = 13.85473523665755
$ws.Range("D20").Value = 8.388885794208015
$ws.Range("E20").Value = 14.32322034182938
$ws.Range("F20").Value = 40.25248078200058
$ws.Range("G20").Value = 3.693454667687043
$ws.Range("J20").Value = 10.94956217152282
$ws.Range("K20").Value = 20.17353712679531
$ws.Range("L20").Value = 9.687947883391059
$ws.Range("N20").Value = 18.45574152989271
$ws.Range("O20").Value = 30.75994841532183
$ws.Range("C21").Value = 13.93697745046404
$ws.Range("D21").Value = 8.432298039871021
$ws.Range("E21").Value = 14.34145255039013
$ws.Range("F21").Value = 40.23545730751543
$ws.Range("G21").Value = 3.690613298564058
$ws.Range("J21").Value = 10.92898541954206
$ws.Range("K21").Value = 20.64277868315912
$ws.Range("L21").Value = 9.675655623668773
$ws.Range("N21").Value = 18.40626537470196
$ws.Range("O21").Value = 30.69807364515523
$ws.Range("C22").Value = 13.99285755567836
$ws.Range("D22").Value = 8.461406758167099
$ws.Range("E22").Value = 14.35522789085492
$ws.Range("F22").Value = 40.23285655466084
$ws.Range("G22").Value = 3.688825626095949
$ws.Range("J22").Value = 10.91664438309813
$ws.Range("K22").Value = 20.94559042508254
$ws.Range("L22").Value = 9.668339587674648
$ws.Range("N22").Value = 18.37558460956905
$ws.Range("O22").Value = 30.66382421315379
$ws.Range("C23").Value = 13.96284229303565
$ws.Range("D23").Value = 8.445805268255986
$ws.Range("E23").Value = 14.34770643664496
$ws.Range("F23").Value = 40.23346424182025
$ws.Range("G23").Value = 3.689773417624319
$ws.Range("J23").Value = 10.92312982017463
$ws.Range("K23").Value = 20.78435784199534
$ws.Range("L23").Value = 9.672178766933451
$ws.Range("N23").Value = 18.39180831672601
$ws.Range("O23").Value = 30.68153630766781
$ws.Range("C24").Value = 13.85342784451141
$ws.Range("D24").Value = 8.388189060645857
$ws.Range("E24").Value = 14.32295392258409
$ws.Range("F24").Value = 40.25290283899816
$ws.Range("G24").Value = 3.693502283352015
$ws.Range("J24").Value = 10.94991719714838
$ws.Range("K24").Value = 20.16580925518847
$ws.Range("L24").Value = 9.688160940745233
$ws.Range("N24").Value = 18.45657817680704
$ws.Range("O24").Value = 30.76106394792165
$ws.Range("C25").Value = 13.74492803788661
$ws.Range("D25").Value = 8.329373863443818
$ws.Range("E25").Value = 14.30431754333315
$ws.Range("F25").Value = 40.31040966775379
$ws.Range("G25").Value = 3.697824202715741
$ws.Range("J25").Value = 10.98357358479732
$ws.Range("K25").Value = 19.48520480958923
$ws.Range("L25").Value = 9.708495478408794
$ws.Range("N25").Value = 18.53357011474601
$ws.Range("O25").Value = 31.14899517283063
